# Requirements.xlsx edit — "promene posle 1. casa predavanja"
#
# - Row 1 (header) + default row heights shrink slightly (50 -> 49.95)
# - Row 2 ("serijska/sekvencijalna datoteka..." requirement) text replaced by a
#   new "rukovanje informacionim resursima" requirement; row grows taller.
# - Two brand-new sub-requirement rows are inserted right after row 2
#   (1.1 / 1.2), each holding only a column-A note.
# - The remaining requirement rows are re-worded ("da napravimo/omogucimo" ->
#   "napraviti/omoguciti") and re-ordered; their numeric effort estimate
#   (column G) travels with the row, so once the two new rows are inserted the
#   existing G values already line up and only the text needs touching.
# - View state: zoom to 153%, selection moves to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the two new sub-requirement rows (pushes old rows 3-7 to 5-9,
#     and the trailing blank spacer rows 10-11 to 12-13) ---
$ws.Rows("3:4").Insert()

# --- row 1 header: slightly shorter row height ---
$ws.Rows(1).RowHeight = 49.95

# --- row 2: replaced requirement text, taller row ---
$ws.Range("A2").Value = "Neophodno je obezbediti rukovanje informacionim resursima. Pod informacionim resursom se podrazumeva kolekcije svojstava objekata posmatranja sa cime ce se manipulisati."
$ws.Rows(2).RowHeight = 94.8

# --- row 3: new sub-requirement 1.1, column A only (borrow A-column's data
#     style from a still-untouched data row so it matches s="8") ---
$ws.Range("A5").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "1.1 inf res treba da ga kreiramo, izmenimo,  sacuvamo."
$ws.Range("B3:G3").Clear()
$ws.Rows(3).RowHeight = 73.95

# --- row 4: new sub-requirement 1.2, column A only ---
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "1.2 da ga sacuvamo, promenimo ime, obrisemo."
$ws.Range("B4:G4").Clear()
$ws.Rows(4).RowHeight = 73.95

# --- row 5: "kolekcija korisnika" requirement, re-worded + trailing period ---
$ws.Range("A5").Value = "Neophodno je napraviti kolekciju korisnika I kolekcije sa podacima o njima."
$ws.Rows(5).RowHeight = 54

# --- row 6: "dodavanje novog korisnika" requirement, re-worded ---
$ws.Range("A6").Value = "Neophodno je omoguciti dodavanje novog korisnika."
$ws.Rows(6).RowHeight = 49.95

# --- row 7: "obrazac za prijavljivanje" requirement, re-worded ---
$ws.Range("A7").Value = "Neophodno je napraviti  obrazac za prijavljivanje korisnika I rukovanjem pravima pristupa."
$ws.Rows(7).RowHeight = 83.4

# --- row 8: "interaktivna pocetna strana" requirement, unchanged wording ---
$ws.Range("A8").Value = "Neophodno je da napravimo interaktivnu pocetnu stranu sa svim alatima za upravljanje sistemom I rukovanje sa greskama."
$ws.Rows(8).RowHeight = 69

# --- row 9: "profil ulogovanog korisnika" requirement, re-worded ---
$ws.Range("A9").Value = "Neophodno je napraviti  profil sa podacima ulogovanog korisnika sa mogucnoscu promene podataka."
$ws.Rows(9).RowHeight = 57.6

# --- trailing blank spacer rows keep the slightly shorter row height too ---
$ws.Rows(12).RowHeight = 49.95
$ws.Rows(13).RowHeight = 49.95

# --- view state: zoom + new selection ---
$excel.ActiveWindow.Zoom = 153
$ws.Range("D6").Select()
